$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new columns P1 and Q1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header's direct formatting (bold font, borders, centered
# alignment) from O1 onto the two newly-added header cells.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-25 ---
# Swap the values in columns I and K, and in columns M and O.
# Then append two new columns (P, Q) with value 2 in every data row.
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value = $kVal   # I <- old K
    $ws.Cells.Item($r, 11).Value = $iVal  # K <- old I
    $ws.Cells.Item($r, 13).Value = $oVal  # M <- old O
    $ws.Cells.Item($r, 15).Value = $mVal  # O <- old M

    $ws.Cells.Item($r, 16).Value = 2      # column P
    $ws.Cells.Item($r, 17).Value = 2      # column Q
}
